$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week")

# Set D587:D612 to Text format so bsecode strings are preserved as text (not auto-numeric)
$ws.Range("D587:D612").NumberFormat = "@"

$ws.Cells.Item(587, 1).Value = 1
$ws.Cells.Item(587, 2).Value = "OFSS"
$ws.Cells.Item(587, 3).Value = "Oracle Financial Services Software Limited"
$ws.Cells.Item(587, 4).Value = "532466"
$ws.Cells.Item(587, 5).Value = 0.33
$ws.Cells.Item(587, 6).Value = 11696.45
$ws.Cells.Item(587, 7).Value = 72200
$ws.Cells.Item(587, 8).Value = "week"
$ws.Cells.Item(587, 9).Value = "29/11/2024 11:38:51"

$ws.Cells.Item(588, 1).Value = 2
$ws.Cells.Item(588, 2).Value = "COFORGE"
$ws.Cells.Item(588, 3).Value = "Coforge (Niit Tech)"
$ws.Cells.Item(588, 4).Value = "532541"
$ws.Cells.Item(588, 5).Value = 0.29
$ws.Cells.Item(588, 6).Value = 8685.85
$ws.Cells.Item(588, 7).Value = 226972
$ws.Cells.Item(588, 8).Value = "week"
$ws.Cells.Item(588, 9).Value = "29/11/2024 11:38:51"

$ws.Cells.Item(589, 1).Value = 3
$ws.Cells.Item(589, 2).Value = "EICHERMOT"
$ws.Cells.Item(589, 3).Value = "Eicher Motors Limited"
$ws.Cells.Item(589, 4).Value = "505200"
$ws.Cells.Item(589, 5).Value = 0.34
$ws.Cells.Item(589, 6).Value = 4831.85
$ws.Cells.Item(589, 7).Value = 407775
$ws.Cells.Item(589, 8).Value = "week"
$ws.Cells.Item(589, 9).Value = "29/11/2024 11:38:51"

$ws.Cells.Item(590, 1).Value = 4
$ws.Cells.Item(590, 2).Value = "TIINDIA"
$ws.Cells.Item(590, 3).Value = "Tube Investments of India Ltd"
$ws.Cells.Item(590, 4).Value = "540762"
$ws.Cells.Item(590, 5).Value = 0.33
$ws.Cells.Item(590, 6).Value = 3590.55
$ws.Cells.Item(590, 7).Value = 353315
$ws.Cells.Item(590, 8).Value = "week"
$ws.Cells.Item(590, 9).Value = "29/11/2024 11:38:51"

$ws.Cells.Item(591, 1).Value = 5
$ws.Cells.Item(591, 2).Value = "ESCORTS"
$ws.Cells.Item(591, 3).Value = "Escorts Limited"
$ws.Cells.Item(591, 4).Value = "500495"
$ws.Cells.Item(591, 5).Value = 1.15
$ws.Cells.Item(591, 6).Value = 3553.9
$ws.Cells.Item(591, 7).Value = 59622
$ws.Cells.Item(591, 8).Value = "week"
$ws.Cells.Item(591, 9).Value = "29/11/2024 11:38:51"

$ws.Cells.Item(592, 1).Value = 6
$ws.Cells.Item(592, 2).Value = "CYIENT"
$ws.Cells.Item(592, 3).Value = "Cyient Limited"
$ws.Cells.Item(592, 4).Value = "532175"
$ws.Cells.Item(592, 5).Value = -0.23
$ws.Cells.Item(592, 6).Value = 1853.45
$ws.Cells.Item(592, 7).Value = 173685
$ws.Cells.Item(592, 8).Value = "week"
$ws.Cells.Item(592, 9).Value = "29/11/2024 11:38:51"

$ws.Cells.Item(593, 1).Value = 7
$ws.Cells.Item(593, 2).Value = "MFSL"
$ws.Cells.Item(593, 3).Value = "Max Financial Services Limited"
$ws.Cells.Item(593, 4).Value = "500271"
$ws.Cells.Item(593, 5).Value = -0.5600000000000001
$ws.Cells.Item(593, 6).Value = 1133.95
$ws.Cells.Item(593, 7).Value = 966284
$ws.Cells.Item(593, 8).Value = "week"
$ws.Cells.Item(593, 9).Value = "29/11/2024 11:38:51"

$ws.Cells.Item(594, 1).Value = 8
$ws.Cells.Item(594, 2).Value = "TATACHEM"
$ws.Cells.Item(594, 3).Value = "Tata Chemicals Limited"
$ws.Cells.Item(594, 4).Value = "500770"
$ws.Cells.Item(594, 5).Value = 0.72
$ws.Cells.Item(594, 6).Value = 1111.75
$ws.Cells.Item(594, 7).Value = 334181
$ws.Cells.Item(594, 8).Value = "week"
$ws.Cells.Item(594, 9).Value = "29/11/2024 11:38:51"

$ws.Cells.Item(595, 1).Value = 9
$ws.Cells.Item(595, 2).Value = "MAXHEALTH"
$ws.Cells.Item(595, 3).Value = "Max Healthcare Institute Ltd"
$ws.Cells.Item(595, 4).Value = "543220"
$ws.Cells.Item(595, 5).Value = 0.3
$ws.Cells.Item(595, 6).Value = 979.75
$ws.Cells.Item(595, 7).Value = 2495600
$ws.Cells.Item(595, 8).Value = "week"
$ws.Cells.Item(595, 9).Value = "29/11/2024 11:38:51"

$ws.Cells.Item(596, 1).Value = 10
$ws.Cells.Item(596, 2).Value = "SYNGENE"
$ws.Cells.Item(596, 3).Value = "Syngene International Limited"
$ws.Cells.Item(596, 4).Value = "539268"
$ws.Cells.Item(596, 5).Value = 2.61
$ws.Cells.Item(596, 6).Value = 940.8
$ws.Cells.Item(596, 7).Value = 1275651
$ws.Cells.Item(596, 8).Value = "week"
$ws.Cells.Item(596, 9).Value = "29/11/2024 11:38:51"

$ws.Cells.Item(597, 1).Value = 11
$ws.Cells.Item(597, 2).Value = "JSL"
$ws.Cells.Item(597, 3).Value = "Jindal Stainless Limited"
$ws.Cells.Item(597, 4).Value = "532508"
$ws.Cells.Item(597, 5).Value = -1.84
$ws.Cells.Item(597, 6).Value = 683.2
$ws.Cells.Item(597, 7).Value = 691415
$ws.Cells.Item(597, 8).Value = "week"
$ws.Cells.Item(597, 9).Value = "29/11/2024 11:38:51"

$ws.Cells.Item(598, 1).Value = 12
$ws.Cells.Item(598, 2).Value = "SONACOMS"
$ws.Cells.Item(598, 3).Value = "Sona BLW Precision Forgings Ltd"
$ws.Cells.Item(598, 4).Value = "543300"
$ws.Cells.Item(598, 5).Value = 1.04
$ws.Cells.Item(598, 6).Value = 671.6
$ws.Cells.Item(598, 7).Value = 1776675
$ws.Cells.Item(598, 8).Value = "week"
$ws.Cells.Item(598, 9).Value = "29/11/2024 11:38:51"

$ws.Cells.Item(599, 1).Value = 13
$ws.Cells.Item(599, 2).Value = "HDFCLIFE"
$ws.Cells.Item(599, 3).Value = "HDFC Life Insurance Company Ltd"
$ws.Cells.Item(599, 4).Value = "540777"
$ws.Cells.Item(599, 5).Value = 0.01
$ws.Cells.Item(599, 6).Value = 657.75
$ws.Cells.Item(599, 7).Value = 5934407
$ws.Cells.Item(599, 8).Value = "week"
$ws.Cells.Item(599, 9).Value = "29/11/2024 11:38:51"

$ws.Cells.Item(600, 1).Value = 14
$ws.Cells.Item(600, 2).Value = "LICHSGFIN"
$ws.Cells.Item(600, 3).Value = "Lic Housing Finance Limited"
$ws.Cells.Item(600, 4).Value = "500253"
$ws.Cells.Item(600, 5).Value = 0.65
$ws.Cells.Item(600, 6).Value = 638.8
$ws.Cells.Item(600, 7).Value = 831816
$ws.Cells.Item(600, 8).Value = "week"
$ws.Cells.Item(600, 9).Value = "29/11/2024 11:38:51"

$ws.Cells.Item(601, 1).Value = 15
$ws.Cells.Item(601, 2).Value = "INDIANB"
$ws.Cells.Item(601, 3).Value = "Indian Bank"
$ws.Cells.Item(601, 4).Value = "532814"
$ws.Cells.Item(601, 5).Value = 0.22
$ws.Cells.Item(601, 6).Value = 574.3
$ws.Cells.Item(601, 7).Value = 2166791
$ws.Cells.Item(601, 8).Value = "week"
$ws.Cells.Item(601, 9).Value = "29/11/2024 11:38:51"

$ws.Cells.Item(602, 1).Value = 16
$ws.Cells.Item(602, 2).Value = "PFC"
$ws.Cells.Item(602, 3).Value = "Power Finance Corporation Limited"
$ws.Cells.Item(602, 4).Value = "532810"
$ws.Cells.Item(602, 5).Value = 0.26
$ws.Cells.Item(602, 6).Value = 495.3
$ws.Cells.Item(602, 7).Value = 7811375
$ws.Cells.Item(602, 8).Value = "week"
$ws.Cells.Item(602, 9).Value = "29/11/2024 11:38:51"

$ws.Cells.Item(603, 1).Value = 17
$ws.Cells.Item(603, 2).Value = "BEL"
$ws.Cells.Item(603, 3).Value = "Bharat Electronics Limited"
$ws.Cells.Item(603, 4).Value = "500049"
$ws.Cells.Item(603, 5).Value = 0.74
$ws.Cells.Item(603, 6).Value = 308
$ws.Cells.Item(603, 7).Value = 23241947
$ws.Cells.Item(603, 8).Value = "week"
$ws.Cells.Item(603, 9).Value = "29/11/2024 11:38:51"

$ws.Cells.Item(604, 1).Value = 18
$ws.Cells.Item(604, 2).Value = "BHEL"
$ws.Cells.Item(604, 3).Value = "Bharat Heavy Electricals Limited"
$ws.Cells.Item(604, 4).Value = "500103"
$ws.Cells.Item(604, 5).Value = -0.55
$ws.Cells.Item(604, 6).Value = 251.09
$ws.Cells.Item(604, 7).Value = 8533794
$ws.Cells.Item(604, 8).Value = "week"
$ws.Cells.Item(604, 9).Value = "29/11/2024 11:38:51"

$ws.Cells.Item(605, 1).Value = 19
$ws.Cells.Item(605, 2).Value = "CUB"
$ws.Cells.Item(605, 3).Value = "City Union Bank Limited"
$ws.Cells.Item(605, 4).Value = "532210"
$ws.Cells.Item(605, 5).Value = -0.08
$ws.Cells.Item(605, 6).Value = 179.53
$ws.Cells.Item(605, 7).Value = 1420556
$ws.Cells.Item(605, 8).Value = "week"
$ws.Cells.Item(605, 9).Value = "29/11/2024 11:38:51"

$ws.Cells.Item(606, 1).Value = 20
$ws.Cells.Item(606, 2).Value = "MANAPPURAM"
$ws.Cells.Item(606, 3).Value = "Manappuram Finance Limited"
$ws.Cells.Item(606, 4).Value = "531213"
$ws.Cells.Item(606, 5).Value = -0.29
$ws.Cells.Item(606, 6).Value = 156.26
$ws.Cells.Item(606, 7).Value = 3356222
$ws.Cells.Item(606, 8).Value = "week"
$ws.Cells.Item(606, 9).Value = "29/11/2024 11:38:51"

$ws.Cells.Item(607, 1).Value = 21
$ws.Cells.Item(607, 2).Value = "IRFC"
$ws.Cells.Item(607, 3).Value = "Indian Railway Finance Corporation Ltd"
$ws.Cells.Item(607, 4).Value = "543257"
$ws.Cells.Item(607, 5).Value = -2.65
$ws.Cells.Item(607, 6).Value = 149.34
$ws.Cells.Item(607, 7).Value = 18576066
$ws.Cells.Item(607, 8).Value = "week"
$ws.Cells.Item(607, 9).Value = "29/11/2024 11:38:51"

$ws.Cells.Item(608, 1).Value = 22
$ws.Cells.Item(608, 2).Value = "HFCL"
$ws.Cells.Item(608, 3).Value = "Himachal Futuristic Communications Limited"
$ws.Cells.Item(608, 4).Value = "500183"
$ws.Cells.Item(608, 5).Value = -3.08
$ws.Cells.Item(608, 6).Value = 129.09
$ws.Cells.Item(608, 7).Value = 20074519
$ws.Cells.Item(608, 8).Value = "week"
$ws.Cells.Item(608, 9).Value = "29/11/2024 11:38:51"

$ws.Cells.Item(609, 1).Value = 23
$ws.Cells.Item(609, 2).Value = "UNIONBANK"
$ws.Cells.Item(609, 3).Value = "Union Bank Of India"
$ws.Cells.Item(609, 4).Value = "532477"
$ws.Cells.Item(609, 5).Value = -0.06
$ws.Cells.Item(609, 6).Value = 121.62
$ws.Cells.Item(609, 7).Value = 23869952
$ws.Cells.Item(609, 8).Value = "week"
$ws.Cells.Item(609, 9).Value = "29/11/2024 11:38:51"

$ws.Cells.Item(610, 1).Value = 24
$ws.Cells.Item(610, 2).Value = "BANKINDIA"
$ws.Cells.Item(610, 3).Value = "Bank Of India"
$ws.Cells.Item(610, 4).Value = "532149"
$ws.Cells.Item(610, 5).Value = -1.25
$ws.Cells.Item(610, 6).Value = 110.5
$ws.Cells.Item(610, 7).Value = 13030139
$ws.Cells.Item(610, 8).Value = "week"
$ws.Cells.Item(610, 9).Value = "29/11/2024 11:38:51"

$ws.Cells.Item(611, 1).Value = 25
$ws.Cells.Item(611, 2).Value = "PNB"
$ws.Cells.Item(611, 3).Value = "Punjab National Bank"
$ws.Cells.Item(611, 4).Value = "532461"
$ws.Cells.Item(611, 5).Value = -1.31
$ws.Cells.Item(611, 6).Value = 104.9
$ws.Cells.Item(611, 7).Value = 30164871
$ws.Cells.Item(611, 8).Value = "week"
$ws.Cells.Item(611, 9).Value = "29/11/2024 11:38:51"

$ws.Cells.Item(612, 1).Value = 26
$ws.Cells.Item(612, 2).Value = "NHPC"
$ws.Cells.Item(612, 3).Value = "Nhpc Limited"
$ws.Cells.Item(612, 4).Value = "533098"
$ws.Cells.Item(612, 5).Value = -2.48
$ws.Cells.Item(612, 6).Value = 81.44
$ws.Cells.Item(612, 7).Value = 25198619
$ws.Cells.Item(612, 8).Value = "week"
$ws.Cells.Item(612, 9).Value = "29/11/2024 11:38:51"

# Reset style for D column so it does not retain an explicit non-default style index
$ws.Range("D587:D612").Style = "Normal"

# Fix day-sheet bsecode cells (D988:D992) from text to numeric
$dayws = $wb.Worksheets.Item("day")
$dayws.Range("D988").Value = 500387
$dayws.Range("D989").Value = 532466
$dayws.Range("D990").Value = 500420
$dayws.Range("D991").Value = 542650
$dayws.Range("D992").Value = 532321